$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:F8 (web planform bug correction)
$data = @(
    @(0.0015, 0.0015, 0,      0.153099998831749,  0),
    @(0.0027, 0.0027, 0,      0.1934999972581863, 0),
    @(0.001,  0.001,  0,      0.1684000045061111, 0),
    @(0.0474, 0.0474, 0.0104, 0.2856999933719635, 0),
    @(0.0234, 0.0234, 0,      0.256199985742569,  0),
    @(0.0021, 0.0021, 0,      0.1269000023603439, 0),
    @(0.0005, 0.0005, 0,      0.07069999724626541,0)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $row++
}
